$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 30
$ws1.Cells.Item(3, 6).Value = 634
$ws1.Cells.Item(6, 6).Value = 5607
$ws1.Cells.Item(7, 6).Value = 1571
$ws1.Cells.Item(8, 6).Value = 164
$ws1.Cells.Item(9, 6).Value = 3190
$ws1.Cells.Item(12, 6).Value = 1325
$ws1.Cells.Item(13, 6).Value = 4448
$ws1.Cells.Item(14, 6).Value = 1061
$ws1.Cells.Item(22, 6).Value = 1006
$ws1.Cells.Item(25, 6).Value = 13
$ws1.Cells.Item(27, 6).Value = 209
$ws1.Cells.Item(28, 6).Value = 3
$ws1.Cells.Item(29, 6).Value = 1108
$ws1.Cells.Item(30, 6).Value = 397
$ws1.Cells.Item(31, 6).Value = 73
$ws1.Cells.Item(32, 6).Value = 194
$ws1.Cells.Item(33, 6).Value = 342
$ws1.Cells.Item(34, 6).Value = 89
$ws1.Cells.Item(37, 6).Value = 2225
$ws1.Cells.Item(38, 6).Value = 1040
$ws1.Cells.Item(42, 6).Value = 344
$ws1.Cells.Item(44, 6).Value = 665
$ws1.Cells.Item(45, 6).Value = 22
$ws1.Cells.Item(46, 6).Value = 424
$ws1.Cells.Item(47, 6).Value = 368
$ws1.Cells.Item(48, 6).Value = 224

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 4

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 771

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 771
$ws4.Cells.Item(3, 6).Value = 30
$ws4.Cells.Item(5, 6).Value = 5607
$ws4.Cells.Item(6, 6).Value = 1571
$ws4.Cells.Item(7, 6).Value = 164
$ws4.Cells.Item(9, 6).Value = 3190
$ws4.Cells.Item(10, 6).Value = 1325
$ws4.Cells.Item(11, 6).Value = 4448
$ws4.Cells.Item(12, 6).Value = 1061
$ws4.Cells.Item(23, 6).Value = 1006
$ws4.Cells.Item(26, 6).Value = 13
$ws4.Cells.Item(28, 6).Value = 209
$ws4.Cells.Item(29, 6).Value = 3
$ws4.Cells.Item(30, 6).Value = 1108
$ws4.Cells.Item(31, 6).Value = 397
$ws4.Cells.Item(32, 6).Value = 194
$ws4.Cells.Item(33, 6).Value = 89
$ws4.Cells.Item(37, 6).Value = 1040
$ws4.Cells.Item(42, 6).Value = 344
$ws4.Cells.Item(43, 6).Value = 665
$ws4.Cells.Item(44, 6).Value = 424
$ws4.Cells.Item(45, 6).Value = 368
$ws4.Cells.Item(46, 6).Value = 224
